$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append the 2025-08-21 portfolio row (A6:D6).
# The date column is stored as plain text (matching the existing rows), so
# the leading apostrophe forces text entry instead of Excel's automatic
# date parsing; resetting the style back to "Normal" afterwards clears the
# quote-prefix formatting flag that the text entry would otherwise leave
# behind, keeping the new row's formatting consistent with rows 2-5.
$ws.Cells.Item(6, 1).Value = "'2025-08-21"
$ws.Cells.Item(6, 1).Style = "Normal"
$ws.Cells.Item(6, 2).Value = 58.15999984741211
$ws.Cells.Item(6, 3).Value = 685.4000244140625
$ws.Cells.Item(6, 4).Value = 321.7999877929688
